$wb = $excel.ActiveWorkbook

# --- Metadata sheet updates ---
$meta = $wb.Worksheets.Item("Metadata")

$meta.Range("B2").Value = "http://linuxforhealth.org/fhir/cdm/StructureDefinition/reference-sequence"
$meta.Range("B3").Value = "8.0.0"
$meta.Range("B8").Value = "2022-11-10T16:00:46+00:00"
$meta.Range("B9").Value = "LinuxForHealth Team"

# --- Elements sheet updates ---
$elements = $wb.Worksheets.Item("Elements")

# Clear the Constraint(s) value for the root Extension row (row 2, column AI)
$elements.Range("AI2").Value = ""

# Keep the Fixed Value (Extension.url row) in sync with the new canonical URL
$elements.Range("Q5").Value = "http://linuxforhealth.org/fhir/cdm/StructureDefinition/reference-sequence"
